# Apply the "Add files via upload" update to the ranking workbook.
# This appends 39 new ranking entries (names/points) to the bottom of the
# tracking sheet, bumps the "update" value in C2, and moves the active
# selection/scroll position to reflect the newly added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# The running "update" total in C2 increased from 26.1 to 28.1
$ws.Cells.Item(2, 3).Value = 28.1

# New ranking rows appended starting at row 318 (name, points)
$newRows = @(
    @("הגר אגמון", 1),
    @("תומר ששון", 1),
    @("אן מרש", 1),
    @("יולי יערי תליו", 1),
    @("ליהי בראל", 1),
    @("תומר ששון", 6),
    @("ליהי בראל", 6),
    @("רומי הרשקוביץ", 1),
    @("עדן ורד מרי", 1),
    @("איתי הראל", 1),
    @("אורי שטרנברג", 1),
    @("תאיו ורד", 1),
    @("איתי בסטקר", 1),
    @("מעיין סטרוזר", 1),
    @("ליאם דיין", 1),
    @("שלו דיין", 1),
    @("הילה שולויס", 1),
    @("יהלי דוייב", 1),
    @("הילה שולויס", 6),
    @("מעיין סטרוזר", 6),
    @("יהלי דוייב", 1),
    @("יולי קזמה", 1),
    @("תומר ששון", 1),
    @("שלו דיין", 1),
    @("תומר ששון", 6),
    @("תומר ששון", 6),
    @("דן פימה", 1),
    @("אביב ואסקז", 1),
    @("ליהי בראל", 1),
    @("ירון גלפנד", 1),
    @("תאיו ורד", 1),
    @("אורי שטרנברג", 1),
    @("יולי יערי תליו", 1),
    @("אן מרש", 1),
    @("קרן רינת פביאן", 1),
    @("ליאם דיין", 1),
    @("איתי הראל", 1),
    @("ליהי בראל", 6),
    @("אן מרש", 6)
)

$startRow = 318
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}

# Move selection / scroll position to the newly added region
$excel.ActiveWindow.ScrollRow = 317
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A337").Select()
